# Script: apply betexplorer.com scrape refresh to 2023 Portugal Liga 3 sheet.
# Two existing match rows had their fixture data (columns F:V — the actual
# match/odds/url fields) re-ordered (column A "Indice" and the B:E
# pais/torneio/temporada/data_partida fields stay anchored to the row), and
# one brand-new fixture row (113, Amora vs Academica) is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 83-88 -----------------------------------------------
# Row 88's fixture data (Atletico CP vs Academica) moves up to row 83;
# rows 83-87's fixture data each shift down by one row.
$block1Last  = $ws.Range("F88:V88").Value2
$block1Shift = $ws.Range("F83:V87").Value2
$ws.Range("F84:V88").Value = $block1Shift
$ws.Range("F83:V83").Value = $block1Last

# --- Block 2: rows 105-108 ----------------------------------------------
# Row 108's fixture data (Lusitania FC vs Braga B) moves up to row 105;
# rows 105-107's fixture data each shift down by one row.
$block2Last  = $ws.Range("F108:V108").Value2
$block2Shift = $ws.Range("F105:V107").Value2
$ws.Range("F106:V108").Value = $block2Shift
$ws.Range("F105:V105").Value = $block2Last

# --- New row 113: Amora vs Academica ------------------------------------
# Clone formatting from the last existing data row (112) so the new row
# picks up the same styles (bold/bordered index cell, date-formatted
# data_partida cell) used throughout the sheet.
$ws.Range("A112:V112").Copy()
$ws.Range("A113:V113").PasteSpecial(-4122)

$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "portugal"
$ws.Range("C113").Value = "liga-3"
$ws.Range("D113").Value = "2023-2024"
$ws.Range("E113").Value = 45262.65625
$ws.Range("F113").Value = "Amora"
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = "Academica"
$ws.Range("I113").Value = 2
$ws.Range("J113").Value = 3.72
$ws.Range("K113").Value = "26/11/2023 18:13"
$ws.Range("L113").Value = 3.39
$ws.Range("M113").Value = "02/12/2023 15:17"
$ws.Range("N113").Value = 3.45
$ws.Range("O113").Value = "26/11/2023 18:13"
$ws.Range("P113").Value = 3.66
$ws.Range("Q113").Value = "02/12/2023 15:17"
$ws.Range("R113").Value = 1.94
$ws.Range("S113").Value = "26/11/2023 18:13"
$ws.Range("T113").Value = 2.1
$ws.Range("U113").Value = "02/12/2023 15:17"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/portugal/liga-3/amora-academica/M7K16K2P/"
